$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price (D) column values are stored as text so numeric-looking
# strings like "1.000" are not coerced into numbers.
$priceCells = @("D2", "D3", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D17", "D19", "D20", "D21", "D23", "D25", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D41", "D42", "D43", "D44", "D45", "D46", "D48", "D49", "D50", "D51")
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '27.287.74'
$ws.Range('E2').Value = '  +0.24%  '
$ws.Range('D3').Value = '1.910.70'
$ws.Range('E3').Value = '  +0.29%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '307.65'
$ws.Range('E5').Value = '  -0.05%  '
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  -0.12%  '
$ws.Range('D7').Value = '0.5266'
$ws.Range('E7').Value = '  +1.21%  '
$ws.Range('D8').Value = '0.3813'
$ws.Range('E8').Value = '  +1.15%  '
$ws.Range('D9').Value = '0.07298'
$ws.Range('E9').Value = '  +0.29%  '
$ws.Range('D10').Value = '22.12'
$ws.Range('E10').Value = '  +4.33%  '
$ws.Range('D11').Value = '0.9016'
$ws.Range('E11').Value = '  -0.29%  '
$ws.Range('D12').Value = '0.08166'
$ws.Range('E12').Value = '  -2.95%  '
$ws.Range('D13').Value = '96.08'
$ws.Range('E13').Value = '  -0.92%  '
$ws.Range('D14').Value = '5.363'
$ws.Range('E14').Value = '  +1.15%  '
$ws.Range('D15').Value = '1.421.52'
$ws.Range('E15').Value = '  -25.32%  '
$ws.Range('E16').Value = '  -0.17%  '
$ws.Range('D17').Value = '0.000008671'
$ws.Range('E17').Value = '  +0.16%  '
$ws.Range('E18').Value = '  +1.61%  '
$ws.Range('D19').Value = '1.000'
$ws.Range('E19').Value = '  -0.11%  '
$ws.Range('D20').Value = '27.312.15'
$ws.Range('E20').Value = '  +0.20%  '
$ws.Range('D21').Value = '5.103'
$ws.Range('E21').Value = '  +0.20%  '
$ws.Range('E22').Value = '  +1.59%  '
$ws.Range('D23').Value = '6.523'
$ws.Range('E23').Value = '  +1.13%  '
$ws.Range('E24').Value = '  +2.21%  '
$ws.Range('D25').Value = '2.307'
$ws.Range('E25').Value = '  -0.56%  '
$ws.Range('D26').Value = '18.26'
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('E27').Value = '  -0.67%  '
$ws.Range('D28').Value = '116.84'
$ws.Range('E28').Value = '  +1.43%  '
$ws.Range('D29').Value = '4.859'
$ws.Range('D30').Value = '4.843'
$ws.Range('E30').Value = '  -1.23%  '
$ws.Range('D31').Value = '0.09261'
$ws.Range('E31').Value = '  -0.14%  '
$ws.Range('D32').Value = '0.8304'
$ws.Range('E32').Value = '  +3.89%  '
$ws.Range('D33').Value = '0.05073'
$ws.Range('E33').Value = '  -0.11%  '
$ws.Range('D34').Value = '1.233'
$ws.Range('E34').Value = '  -0.91%  '
$ws.Range('D35').Value = '2.996'
$ws.Range('E35').Value = '  +1.08%  '
$ws.Range('D36').Value = '3.365'
$ws.Range('E36').Value = '  -1.90%  '
$ws.Range('E37').Value = '  +4.74%  '
$ws.Range('D38').Value = '0.5806'
$ws.Range('E38').Value = '  +0.03%  '
$ws.Range('D39').Value = '0.02003'
$ws.Range('E39').Value = '  -0.12%  '
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('D41').Value = '9.258'
$ws.Range('E41').Value = '  +2.11%  '
$ws.Range('D42').Value = '6.582'
$ws.Range('E42').Value = '  -0.37%  '
$ws.Range('D43').Value = '116.67'
$ws.Range('E43').Value = '  -0.04%  '
$ws.Range('D44').Value = '0.1527'
$ws.Range('E44').Value = '  +0.44%  '
$ws.Range('D45').Value = '0.4928'
$ws.Range('E45').Value = '  +1.08%  '
$ws.Range('D46').Value = '10.19'
$ws.Range('E46').Value = '  +0.47%  '
$ws.Range('E47').Value = '  -0.17%  '
$ws.Range('D48').Value = '1.646'
$ws.Range('E48').Value = '  +0.34%  '
$ws.Range('D49').Value = '38.92'
$ws.Range('E49').Value = '  +3.08%  '
$ws.Range('D50').Value = '0.06181'
$ws.Range('E50').Value = '  +3.67%  '
$ws.Range('D51').Value = '64.44'
$ws.Range('E51').Value = '  +0.59%  '
